# "some text added by anand"
#
# Adds a new line of text below the existing table on Sheet1 (new shared
# string "//this text added by anand" written to cell B33), which grows the
# sheet's used range from A1:C30 to A1:C33, then leaves that new cell
# selected/active and scrolled into view - matching the author's final
# on-screen state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new text - this also appends a new entry to the shared strings
# table.
$ws.Range("B33").Value = "//this text added by anand"

# Scroll so row 13 is at the top of the window and select/activate the new
# cell, matching the saved view state.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B33").Select()
